# 00-index (agenda) slide: update the per-topic time allocations.
#
#   3. 全学サービスの全体像 (田浦 15分)  -> (田浦 20分)
#   4. セキュリティと在宅勤務 (玉造 15分)  -> (玉造 20分)
#   5. uteleconとサポータについて (鈴木 20分) -> (鈴木 10分)
#
# Each number lives in its own run inside the single "agenda" content
# placeholder, so the safest edit is a targeted Characters() replacement
# of just the digits rather than reassigning the whole TextRange (which
# would blow away the existing run/formatting boundaries).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Full text (129 chars) is:
#   "1. ...(田浦 10分)\r2. ...(柴山 30分)\r3. ...(田浦 15分)\r4. ...(玉造 15分)\r5. ...(鈴木 20分)"
# -> the three numbers we need to touch start at 1-based offsets 73, 96, 126.

$tr.Characters(73, 2).Text = "20"   # item 3: 田浦 15 -> 20
$tr.Characters(96, 2).Text = "20"   # item 4: 玉造 15 -> 20
$tr.Characters(126, 2).Text = "10"  # item 5: 鈴木 20 -> 10
